$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape 36 ("TextBox 75", id=76): "(9) Send push notification to Firebase" ---
# Split into "(10) " + "Send push notification to Firebase", and nudge the
# shape's rotated bounding box (off/ext) to the new position/size.
$shp36 = $s.Shapes.Item(36)
$tr36 = $shp36.TextFrame.TextRange
$ch36 = $tr36.Characters(1, 4)
$ch36.Text = "(10) "

$shp36.Left = 519.9729633608393
$shp36.Top = 46.69334706662417
$shp36.Width = 100.9254780483774
$shp36.Height = 58.16251968503937

# --- Shape 37 ("TextBox 76", id=77): "(10) Send " / "push notification to iOS app" ---
# Renumber (10)->(11) and split "(" / "11) " / "Send " into separate runs.
$shp37 = $s.Shapes.Item(37)
$tr37 = $shp37.TextFrame.TextRange
$chSend37 = $tr37.Characters(6, 5)
$chSend37.Text = "Send "
$ch11_37 = $tr37.Characters(2, 4)
$ch11_37.Text = "11) "

# --- Shape 38 ("Rectangle 77", id=78): "(9) Collect FCM token from database" ---
# Split into "(9) " + "Collect FCM token from database".
$shp38 = $s.Shapes.Item(38)
$tr38 = $shp38.TextFrame.TextRange
$ch38 = $tr38.Characters(1, 4)
$ch38.Text = "(9) "
